$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row just below the header (row 1), shifting all the
# existing price-history rows down by one.
$ws.Rows.Item(2).Insert()

# The cell below (old row 2, now row 3) still carries the normal
# "General" style used throughout column A. Temporarily mark the new
# cell as Text so that typing a date-shaped string does not get
# auto-converted into a date serial number, then restore its style to
# match the rest of the column so no stray formatting is introduced.
$ws.Cells.Item(2, 1).NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = "2026-02-12"
$ws.Cells.Item(2, 1).Style = $ws.Cells.Item(3, 1).Style

$ws.Cells.Item(2, 2).Value = 783.5
$ws.Cells.Item(2, 3).Value = 1112
$ws.Cells.Item(2, 4).Value = 3610
